$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data in row 10 (student "Кику Станислав")
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 5
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 5
$ws.Range("G10").Value = 5
$ws.Range("H10").Value = 5
$ws.Range("I10").Value = 5

# Update row heights and thick top/bottom borders for rows 9, 10, 11
$ws.Rows("9").RowHeight = 14
$ws.Rows("10").RowHeight = 14
$ws.Rows("11").RowHeight = 14

$ws.Range("A9:K9").Borders.Item(8).LineStyle = 1
$ws.Range("A9:K9").Borders.Item(8).Weight = 4
$ws.Range("A9:K9").Borders.Item(9).LineStyle = 1
$ws.Range("A9:K9").Borders.Item(9).Weight = 4

$ws.Range("A10:K10").Borders.Item(8).LineStyle = 1
$ws.Range("A10:K10").Borders.Item(8).Weight = 4
$ws.Range("A10:K10").Borders.Item(9).LineStyle = 1
$ws.Range("A10:K10").Borders.Item(9).Weight = 4

$ws.Range("A11:K11").Borders.Item(8).LineStyle = 1
$ws.Range("A11:K11").Borders.Item(8).Weight = 4
$ws.Range("A11:K11").Borders.Item(9).LineStyle = 1
$ws.Range("A11:K11").Borders.Item(9).Weight = 4

# Update selection to J3
$ws.Range("J3").Select()
